# Apply "Updated profiles with eliminations" changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark eliminated profiles (column C) with a 1 instead of 0.
$ws.Range("C15").Value = 1
$ws.Range("C27").Value = 1
$ws.Range("C38").Value = 1
$ws.Range("C56").Value = 1
$ws.Range("C58").Value = 1

# Update the view: scroll so row 44 is the top visible row, and move the
# active selection to C56.
$window = $excel.ActiveWindow
$window.ScrollRow = 44
$ws.Range("C56").Select()
